# Weekly fruit/vegetable price update: insert a new record as row 197,
# pushing all existing rows from 197 downward by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197 (shifts rows 197..322 down to 198..323,
# and inherits formatting such as the date-style on column D).
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A197").Value = 4
$ws.Range("B197").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C197").Value = "Los Lagos"
$ws.Range("D197").Value = 44680
$ws.Range("E197").Value = 10
$ws.Range("F197").Value = 100112023
$ws.Range("G197").Value = "Brócoli"
$ws.Range("H197").Value = "Sin especificar"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 700
$ws.Range("K197").Value = 1600
$ws.Range("L197").Value = 1700
$ws.Range("M197").Value = 1650
$ws.Range("N197").Value = "$/unidad"
$ws.Range("O197").Value = "Región Metropolitana"
$ws.Range("P197").Value = 1650
$ws.Range("Q197").Value = 1
$ws.Range("R197").Value = "Hortaliza"
